# Sync the T_SERVER_ZONE schema sheet: mark createdAt/updatedAt as NOT (NULL
# not allowed), matching the other timestamp/id columns, and touch up the
# IsNull cells that already said "NOT" so the sheet is consistent.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("T_SERVER_ZONE")
$ws.Activate()

# column C = "IsNull" -> value "NOT" marks the field as NOT NULL
$ws.Range("C3").Value = "NOT"
$ws.Range("C13").Value = "NOT"
$ws.Range("C14").Value = "NOT"
$ws.Range("C18").Value = "NOT"

# page setup was touched while reviewing/printing the sheet
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# cursor ended on C18 after the edits
$ws.Range("C18").Select() | Out-Null
